$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.388675451278687
$ws.Range("B1").Value = 1.499651312828064
$ws.Range("C1").Value = 1.687732815742493
$ws.Range("D1").Value = 2.611310482025146
$ws.Range("E1").Value = -1
